$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 55, shifting existing rows 55-82 down to 56-83.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new record's data.
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44813
$ws.Cells.Item(55, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = 300000001
$ws.Cells.Item(55, 7).Value = "Rabanito"
$ws.Cells.Item(55, 8).Value = "Sin especificar"
$ws.Cells.Item(55, 9).Value = "Primera"
$ws.Cells.Item(55, 10).Value = 30
$ws.Cells.Item(55, 11).Value = 8000
$ws.Cells.Item(55, 12).Value = 8000
$ws.Cells.Item(55, 13).Value = 8000
$ws.Cells.Item(55, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(55, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(55, 16).Value = 667
$ws.Cells.Item(55, 17).Value = 12
$ws.Cells.Item(55, 18).Value = "Hortaliza"
